# new start (version 0.0.10) after CRAN submission of version 0.0.9
#
# The accrual denominator changed from N=6 to N=3, so the two header cells
# that reference the total accrual count are updated, and every percentage
# cell (columns C and E, rows 10-24) is recomputed against the new
# denominator (count / 3 * 100). Column B/D counts and the category labels
# (column A) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell while preserving that
# cell's existing style. A plain `Range.Value = "33.33"` assignment lets
# Excel's COM layer auto-coerce number-looking strings into real numbers
# (losing the text type / formatting, e.g. the trailing zeros or the
# leading-space variants used in this table). Routing the text through a
# formula ("="33.33"") and pasting *values only* (xlPasteValues) keeps it
# as genuine text without touching number format / style of the target
# cell (unlike forcing text via NumberFormat="@" or an apostrophe prefix,
# both of which mutate the cell's style index).
$helper = $ws.Range("ZZ1")
function Set-TextValue([string]$cellRef, [string]$text) {
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# Header row (row 9): accrual denominator N=6 -> N=3
Set-TextValue "C9" "the % of subjects that this comprises of the total accrual (N=3)"
Set-TextValue "E9" "% of the subjects that this comprises of the total accrual (N=3)"

# Recomputed percentages, row by row (category label + counts unchanged)
Set-TextValue "C10" "100.00"
Set-TextValue "E10" "66.67"

Set-TextValue "C11" " 33.33"

Set-TextValue "C12" "100.00"
Set-TextValue "E12" "33.33"

Set-TextValue "C13" " 66.67"
Set-TextValue "E13" "33.33"

Set-TextValue "C14" " 66.67"
Set-TextValue "E14" "66.67"

Set-TextValue "C15" " 33.33"
Set-TextValue "E15" "33.33"

Set-TextValue "C16" " 33.33"
Set-TextValue "E16" "33.33"

Set-TextValue "C17" "100.00"
Set-TextValue "E17" "66.67"

Set-TextValue "C18" " 66.67"
Set-TextValue "E18" "33.33"

Set-TextValue "C19" " 33.33"

Set-TextValue "C20" " 33.33"

Set-TextValue "C21" " 33.33"

Set-TextValue "C22" " 66.67"

Set-TextValue "C23" " 66.67"

Set-TextValue "C24" " 33.33"

# Remove the scratch helper cell so it doesn't leak into the saved sheet.
$helper.Clear()
